$d = $word.ActiveDocument

# 1. Remove the "License Suspension." paragraph entirely (merging it away), including
#    the paragraph mark that separates it from the preceding "Fines and Costs." paragraph
#    and the paragraph mark that separates it from the following (empty) paragraph.
$d.Range(1409, 1517).Delete() | Out-Null

# 2. Date change: "on March 13, 2022." -> "on March 15, 2022." (arraignment date)
$d.Content.Find.Execute("on March 13, 2022.", $false, $false, $false, $false, $false, $true, 1, $false, "on March 15, 2022.", 2) | Out-Null

# 3. "Defendant was represented by  , Private Counsel." -> "Defendant waived right to counsel."
$d.Content.Find.Execute("Defendant was represented by  , Private Counsel.", $false, $false, $false, $false, $false, $true, 1, $false, "Defendant waived right to counsel.", 2) | Out-Null

# 4. Remove the visible " R.C. 2943.031." sentence that followed the deportation-advisement text
#    (a hidden/vanish duplicate of this text remains untouched further down the paragraph)
$d.Content.Find.Execute("law. R.C. 2943.031. ", $false, $false, $false, $false, $false, $true, 1, $false, "law. ", 2) | Out-Null

# 5. Date change: "March 13, 2022" -> "March 15, 2022" (fines/costs due date)
$d.Content.Find.Execute("March 13, 2022", $false, $false, $false, $false, $false, $true, 1, $false, "March 15, 2022", 2) | Out-Null
